$wb = $excel.ActiveWorkbook

# --- Status text update: "Ready for handoff" -> "In Translation" ---
# This shared text shows up on the Overview sheet (zh-cn / de-de status
# columns) as well as on each per-locale detail sheet's "Status" column.

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2:F4").Value = "In Translation"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2:C4").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2:C4").Value = "In Translation"

# --- Column widths: the status columns shrink to fit the new text ---
# (report regeneration re-autosizes the "zh-cn"/"de-de" status columns)

$ws1.Range("E1").ColumnWidth = 12.576851254417766
$ws1.Range("F1").ColumnWidth = 12.576851254417766
$ws2.Range("C1").ColumnWidth = 12.576851254417766
$ws3.Range("C1").ColumnWidth = 12.576851254417766
